$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "imipramina"
$ws.Range("B2").Value = "-0.14 (-5.51;  5.23)"
$ws.Range("D2").Value = "."
$ws.Range("F2").Value = "4.46 ( 1.77;  7.16)"

$ws.Range("A3").Value = "1.61 ( -3.23;  6.45)"
$ws.Range("B3").Value = "amitriptilina"
$ws.Range("F3").Value = "1.07 (-4.30;  6.44)"

$ws.Range("A4").Value = "2.31 ( -1.32;  5.94)"
$ws.Range("B4").Value = "0.70 ( -4.72;  6.12)"
$ws.Range("C4").Value = "nortriptilina"
$ws.Range("F4").Value = "2.12 (-0.31;  4.55)"

$ws.Range("A5").Value = "3.12 ( -2.90;  9.13)"
$ws.Range("B5").Value = "1.50 ( -5.73;  8.74)"
$ws.Range("C5").Value = "0.81 ( -5.09;  6.70)"
$ws.Range("D5").Value = "fluoxetina"
$ws.Range("F5").Value = "1.32 (-4.06;  6.69)"

$ws.Range("A6").Value = "4.70 ( -1.31; 10.72)"
$ws.Range("B6").Value = "3.09 ( -4.14; 10.33)"
$ws.Range("C6").Value = "2.39 ( -3.51;  8.29)"
$ws.Range("D6").Value = "1.59 ( -6.02;  9.19)"
$ws.Range("E6").Value = "citalopram"
$ws.Range("F6").Value = "-0.27 (-5.65;  5.10)"

$ws.Range("A7").Value = "4.43 (  1.73;  7.13)"
$ws.Range("B7").Value = "2.82 ( -2.02;  7.67)"
$ws.Range("C7").Value = "2.12 ( -0.31;  4.55)"
$ws.Range("D7").Value = "1.32 ( -4.06;  6.69)"
$ws.Range("E7").Value = "-0.27 ( -5.65;  5.10)"
$ws.Range("F7").Value = "veículo"
